# Atualizado por script em 08-11-2023 02:45
# Appends two new match rows (114 and 115) to the betexplorer-style sheet,
# mirroring the formatting of the existing data rows (bold/bordered/centered
# style on column A, datetime-formatted style on column E) and keeping the
# "looks like a number" text columns (D, K, M, O, Q, S, U) stored as text
# rather than letting Excel auto-coerce them to numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 114 -------------------------------------------------------------

# Column A: copy the number-style (bold, centered, bordered) from the last
# existing row, then overwrite with the new index value.
$ws.Cells.Item(113, 1).Copy()
$ws.Cells.Item(114, 1).PasteSpecial(-4122)
$ws.Cells.Item(114, 1).Value = 113

# Column B-D: plain text fields. D looks numeric ("2023"), so pull it in via
# PasteSpecial(values) from a cell that is already text-typed, which keeps
# it text without touching its (default) style.
$ws.Cells.Item(114, 2).Value = "paraguay"
$ws.Cells.Item(114, 3).Value = "primera-division"
$ws.Cells.Item(113, 4).Copy()
$ws.Cells.Item(114, 4).PasteSpecial(-4163)

# Column E: copy the datetime number-format style, then set the serial value.
$ws.Cells.Item(113, 5).Copy()
$ws.Cells.Item(114, 5).PasteSpecial(-4122)
$ws.Cells.Item(114, 5).Value = 45237.91666666666

$ws.Cells.Item(114, 6).Value = "Guairena"
$ws.Cells.Item(114, 7).Value = 2
$ws.Cells.Item(114, 8).Value = "Ameliano"
$ws.Cells.Item(114, 9).Value = 2
$ws.Cells.Item(114, 10).Value = 3.04
$ws.Cells.Item(114, 11).Value = "03/11/2023 22:12"
$ws.Cells.Item(114, 12).Value = 3.54
$ws.Cells.Item(114, 13).Value = "07/11/2023 21:56"
$ws.Cells.Item(114, 14).Value = 3.39
$ws.Cells.Item(114, 15).Value = "03/11/2023 22:12"
$ws.Cells.Item(114, 16).Value = 3.43
$ws.Cells.Item(114, 17).Value = "07/11/2023 21:57"
$ws.Cells.Item(114, 18).Value = 2.38
$ws.Cells.Item(114, 19).Value = "03/11/2023 22:12"
$ws.Cells.Item(114, 20).Value = 2.17
$ws.Cells.Item(114, 21).Value = "07/11/2023 21:57"
$ws.Cells.Item(114, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/guairena-fc-sportivo-ameliano/YyLHbgcJ/"

# ---- Row 115 -------------------------------------------------------------

$ws.Cells.Item(114, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4122)
$ws.Cells.Item(115, 1).Value = 114

$ws.Cells.Item(115, 2).Value = "paraguay"
$ws.Cells.Item(115, 3).Value = "primera-division"
$ws.Cells.Item(114, 4).Copy()
$ws.Cells.Item(115, 4).PasteSpecial(-4163)

$ws.Cells.Item(114, 5).Copy()
$ws.Cells.Item(115, 5).PasteSpecial(-4122)
$ws.Cells.Item(115, 5).Value = 45238.02083333334

$ws.Cells.Item(115, 6).Value = "Olimpia Asuncion"
$ws.Cells.Item(115, 7).Value = 2
$ws.Cells.Item(115, 8).Value = "Sp. Luqueno"
$ws.Cells.Item(115, 9).Value = 1
$ws.Cells.Item(115, 10).Value = 1.81
$ws.Cells.Item(115, 11).Value = "04/11/2023 00:42"
$ws.Cells.Item(115, 12).Value = 1.96
$ws.Cells.Item(115, 13).Value = "08/11/2023 00:21"
$ws.Cells.Item(115, 14).Value = 3.93
$ws.Cells.Item(115, 15).Value = "04/11/2023 00:42"
$ws.Cells.Item(115, 16).Value = 3.47
$ws.Cells.Item(115, 17).Value = "08/11/2023 00:29"
$ws.Cells.Item(115, 18).Value = 3.91
$ws.Cells.Item(115, 19).Value = "04/11/2023 00:42"
$ws.Cells.Item(115, 20).Value = 4.24
$ws.Cells.Item(115, 21).Value = "08/11/2023 00:29"
$ws.Cells.Item(115, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/olimpia-asuncion-sp-luqueno/nu9aNxzK/"
